$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Anmolpreet Singh"

# Insert a new column A, shifting existing data right
$ws.Range("A1").EntireColumn.Insert()

# Populate new column A
$ws.Range("A1").Value = "matchNo"
$ws.Range("A2").Value = "30th"
